$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.717.44"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "3.303.23"
$ws.Range("E3").Value = "  -5.61%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "181.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "530.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.605"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "3.298.22"
$ws.Range("E8").Value = "  -5.55%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.618"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.134"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.49%  "
$ws.Range("D15").Value = "3.831.55"
$ws.Range("E15").Value = "  -5.68%  "
$ws.Range("D16").Value = "3.304.56"
$ws.Range("E16").Value = "  -5.71%  "
$ws.Range("E17").Value = "  -4.75%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.39%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "64.566.02"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.962"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "652.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.38%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.397"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "0.0₃0702"
$ws.Range("E41").Value = "  +4.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.127"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").Value = "2.875.54"
$ws.Range("E43").Value = "  -6.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0402"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.26%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.55%  "
